$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SPY 161021C00223000, C)
$ws.Range("E2").Value = 0.16
$ws.Range("F2").Value = 0.17
$ws.Range("G2").Value = 0.16
$ws.Range("K2").Value = 0.0635327980336118
$ws.Range("L2").Value = 0.0210419003076321
$ws.Range("M2").Value = -0.0117312964988799
$ws.Range("N2").Value = 0.0816093863996498
$ws.Range("O2").Value = 42632.9111111111
$ws.Range("Q2").Value = 0.0936522598789154
$ws.Range("S2").Value = 213.6
$ws.Range("T2").Value = -495

# Row 3 (SPY 161021C00228000, C)
$ws.Range("G3").Value = 0.03
$ws.Range("K3").Value = 0.00729229547981208
$ws.Range("L3").Value = 0.00353121961828463
$ws.Range("M3").Value = -0.00180588978926703
$ws.Range("N3").Value = 0.0176913090918269
$ws.Range("O3").Value = 42632.9111111111
$ws.Range("Q3").Value = 0.089888083175986
$ws.Range("S3").Value = 213.6

# Row 4 (SPY 161021P00195000, P)
$ws.Range("E4").Value = 0.44
$ws.Range("F4").Value = 0.45
$ws.Range("G4").Value = 0.47
$ws.Range("K4").Value = -0.0695550996747478
$ws.Range("L4").Value = 0.00989485314167953
$ws.Range("M4").Value = -0.0285708991382583
$ws.Range("N4").Value = 0.100302169534225
$ws.Range("O4").Value = 42632.9111111111
$ws.Range("Q4").Value = 0.215777176360249
$ws.Range("S4").Value = 213.6
$ws.Range("T4").Value = 1335

# Row 5 (SPY 161021P00200000, P)
$ws.Range("E5").Value = 0.74
$ws.Range("F5").Value = 0.75
$ws.Range("G5").Value = 0.74
$ws.Range("K5").Value = -0.115897941048817
$ws.Range("L5").Value = 0.0161839028817106
$ws.Range("M5").Value = -0.0372310404313283
$ws.Range("N5").Value = 0.135854956008718
$ws.Range("O5").Value = 42632.9111111111
$ws.Range("Q5").Value = 0.192828773389221
$ws.Range("S5").Value = 213.58
$ws.Range("T5").Value = -2235
